# 📊 Weekly driver report update for 2025-04-28
#
# - Updates the single "Bad Drivers" row with this week's figures.
# - A newly observed driver version (21.40.1.3) enters the "Good Drivers"
#   table at the top; the existing good-driver rows shift down by one,
#   and their "Total Samples" counters are refreshed with this week's
#   cumulative totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Bad Drivers (row 3) + totals (row 4)
# ---------------------------------------------------------------------
$ws.Cells.Item(3,1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.1.2"
$ws.Cells.Item(3,2).Value = 2
$ws.Cells.Item(3,3).Value = 214
$ws.Cells.Item(3,4).Value = 98.59999999999999

$ws.Cells.Item(4,2).Value = 2
$ws.Cells.Item(4,3).Value = 214

# ---------------------------------------------------------------------
# Good Drivers table: make room for the new row at the top by shifting
# rows 12-17 down to 13-18. Copy bottom-up so sources aren't clobbered
# before they are read. Columns A, C, D, E move as-is; column B (Total
# Samples) gets this week's refreshed counter for each driver.
# ---------------------------------------------------------------------

# Materialize cell formatting on the newly-used row 18 (it was blank
# before, so copy the style pattern down from row 17 first).
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)

# Seed C18 with a placeholder so it keeps its (empty) text type once the
# real (also-empty) value is pasted over it below - a truly blank cell
# would otherwise vanish entirely instead of staying an empty text cell.
$ws.Cells.Item(18,3).Value = "placeholder"

# Row 17 -> Row 18 ("...21.60.2.1"); Total Samples unchanged (56018)
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E18").PasteSpecial(-4163)
$ws.Cells.Item(18,2).Value = 56018

# Row 16 -> Row 17 ("...21.70.0.6"); Total Samples 113652 -> 117653
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4163)
$ws.Cells.Item(17,2).Value = 117653

# Row 15 -> Row 16 ("...21.110.3.2"); Total Samples 59673 -> 65425
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4163)
$ws.Cells.Item(16,2).Value = 65425

# Row 14 -> Row 15 ("...22.50.1.1"); Total Samples 34244 -> 35355
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4163)
$ws.Cells.Item(15,2).Value = 35355

# Row 13 -> Row 14 ("...22.80.0.9"); Total Samples 77999 -> 79953
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4163)
$ws.Cells.Item(14,2).Value = 79953

# Row 12 -> Row 13 ("...23.100.0.4"); Total Samples 449371 -> 486214
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4163)
$ws.Cells.Item(13,2).Value = 486214

# New driver inserted at row 12 ("...21.40.1.3"); no vintage date yet,
# so the driver-vintage cell (E12) is left blank.
$ws.Cells.Item(12,1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Cells.Item(12,2).Value = 11128
$ws.Cells.Item(12,4).Value = 100
$ws.Cells.Item(12,5).ClearContents()

$excel.CutCopyMode = 0
